$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = "55"
$ws.Range("H6").Value = "65"
$ws.Range("I6").Value = "64"
$ws.Range("G5").Value = "53"
$ws.Range("J5").Value = "166"
$ws.Range("K5").Value = "169"
$ws.Range("L5").Value = "335"
$ws.Range("G2").Value = "12"
$ws.Range("H5").Value = "58"

$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("G11").Select()
